$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "114×9=1026" "971×7=6797"
Replace-Text "345×6=2070" "536×6=3216"
Replace-Text "725×7=5075" "433×8=3464"
Replace-Text "961×3=2883" "257×8=2056"
Replace-Text "205×8=1640" "985×5=4925"
Replace-Text "824×2=1648" "676×8=5408"
Replace-Text "217×2=434" "329×2=658"
Replace-Text "343×8=2744" "347×6=2082"
Replace-Text "810×5=4050" "197×8=1576"
Replace-Text "640×7=4480" "918×9=8262"
Replace-Text "285×4=1140" "291×7=2037"
Replace-Text "188×5=940" "346×2=692"
Replace-Text "412×3=1236" "475×9=4275"
Replace-Text "986×8=7888" "380×3=1140"
Replace-Text "683×5=3415" "310×5=1550"
Replace-Text "682×6=4092" "713×5=3565"
Replace-Text "361×2=722" "773×4=3092"
Replace-Text "543×2=1086" "408×6=2448"
Replace-Text "820×9=7380" "242×4=968"
Replace-Text "691×9=6219" "574×4=2296"
Replace-Text "822×4=3288" "433×8=3464"
Replace-Text "838×3=2514" "370×9=3330"
Replace-Text "163×3=489" "436×5=2180"
Replace-Text "969×2=1938" "136×6=816"
Replace-Text "991×6=5946" "333×9=2997"

Write-Output "Done replacing $($d.Content.Text.Length) chars of content"
